# Barcodes_Test.xlsx tutorial fix-up: trim the fake-flowcell adapter
# sequences from 10 characters down to 8 (A/C/G/T homopolymer "barcodes")
# and nudge the workbook's tab-bar split ratio.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Window chrome: widen the sheet-tab area relative to the horizontal
# scrollbar (stored as tabRatio in the saved workbook; Excel's TabRatio
# property is a 0..1 fraction, so 987/1000 -> 0.987).
$excel.ActiveWindow.TabRatio = 0.987

# Shorten the four sequence values in column C (rows 9-12) from 10 to 8
# repeated bases, matching the shortened fake barcodes used later in the
# tutorial.
$ws.Range("C9").Value = "AAAAAAAA"
$ws.Range("C10").Value = "CCCCCCCC"
$ws.Range("C11").Value = "GGGGGGGG"
$ws.Range("C12").Value = "TTTTTTTT"
